# "Avanços feitos em Aula" - trims / adds a few requirement bullets in the
# "Padaria Doce Sabor" case study section.

$d = $word.ActiveDocument

# 0) Intro paragraph - split the sentence "...finalizando a compra." from
#    "Na leitura da comanda..." (the sentence text itself is unchanged;
#    Word just re-marks the cursor's last-edit position here, moving the
#    auto-managed "_GoBack" bookmark to sit between the two sentences).
$splitPoint = $d.Content.Duplicate
$splitPoint.Find.Execute("finalizando a compra. ") | Out-Null
$splitPoint.Collapse(0)
$d.Bookmarks.Add("_GoBack", $splitPoint) | Out-Null

# 1) RF: paragraph - remove the now-redundant ", calcular valores gastos"
#    clause right before "e geração de relatórios."
$find = $d.Content.Find
$find.Execute(
    ", calcular valores gastos e geração de relatórios", $true, $false, $false,
    $false, $false, $true, 1, $false,
    " e geração de relatórios", 2) | Out-Null

# 2) RNF: paragraph - drop "permitir saída do cliente após pagamento dos
#    produtos consumidos, " (kept the rest of the sentence intact).
$find = $d.Content.Find
$find.Execute(
    "status da NFE-C, permitir saída do cliente após pagamento dos produtos consumidos, carregar entrada de estoque por XML",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "status da NFE-C, carregar entrada de estoque por XML", 2) | Out-Null

# 3) Domínio: paragraph - append new domain rules after the stock-control
#    sentence.
$find = $d.Content.Find
$find.Execute(
    "não deixar o estoque ficar negativo;", $true, $false, $false, $false,
    $false, $true, 1, $false,
    "não deixar o estoque ficar negativo; cálculo das comandas; condição para saída do estabelecimento;", 2) | Out-Null

# 4) Revistas em quadrinhos RF paragraph - remove "agenda de datas dos
#    empréstimos; " clause.
$find = $d.Content.Find
$find.Execute(
    "pegaram revistas emprestadas; agenda de datas dos empréstimos; ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "pegaram revistas emprestadas; ", 2) | Out-Null
